$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks so we can rebuild them cleanly in final row order
$ws.Cells.Hyperlinks.Delete()

# Column D width: 30 -> 41. ColumnWidth round-trips through Excels char-width
# formula and always adds back ~0.8333, so feed it 40.1666666666667 to land on
# a stored width of exactly 41 (matches target dimension/col width in the diff).
$ws.Columns.Item(4).ColumnWidth = 40.16666666666666

# Row 2
$ws.Range("A2").Value = '2025-12-01 01:44:37'
$ws.Range("B2").Value = '【AI系勉強会】「Google Gravity」開発事例発表者募集!個人開発をプレゼンしませんか?'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5443957'
$ws.Range("G2").Value = 360
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

# Row 3
$ws.Range("A3").Value = '2025-12-01 01:44:37'
$ws.Range("B3").Value = '【急募】BlockChainとSolidityに精通したプログラマー募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '80,000 円 ~ 90,000 円 / 募集期間 1 日、取引期間 0 日'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5443998'
$ws.Range("G3").Value = 298
$ws.Range("H3").Value = '🔥AI,Ai'

# Row 4
$ws.Range("A4").Value = '2025-12-01 01:44:37'
$ws.Range("B4").Value = '顧客のBtoB向けの管理画面サービスの構築とAPI連携'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5444251'
$ws.Range("G4").Value = 220
$ws.Range("H4").Value = '🔥API ◇管理'

# Row 5
$ws.Range("A5").Value = '2025-12-01 01:44:37'
$ws.Range("B5").Value = '【自動化】ニュースサイト情報をX(旧Twitter)へ投稿するシステム開発'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5444198'
$ws.Range("G5").Value = 205
$ws.Range("H5").Value = '◆開発,システム開発 ◇サイト'

# Row 6
$ws.Range("A6").Value = '2025-12-01 01:44:37'
$ws.Range("B6").Value = '外部WEB予約サイト一元管理システム開発|長期保守パートナー募集'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5444378'
$ws.Range("G6").Value = 170
$ws.Range("H6").Value = '◆開発,システム開発 ◇サイト'

# Row 7
$ws.Range("A7").Value = '2025-12-01 01:44:37'
$ws.Range("B7").Value = '【日本人限定/継続案件】Node.jsエンジニア募集(スクレイピング機能開発)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5444489'
$ws.Range("G7").Value = 155
$ws.Range("H7").Value = '◆開発,Node.js'

# Row 8
$ws.Range("A8").Value = '2025-12-01 01:44:37'
$ws.Range("B8").Value = '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5251319'
$ws.Range("G8").Value = 135
$ws.Range("H8").Value = '◆ツール,スクレイピング ◇サイト'

# Row 9
$ws.Range("A9").Value = '2025-12-01 01:44:37'
$ws.Range("B9").Value = '【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5444141'
$ws.Range("G9").Value = 68
$ws.Range("H9").Value = '◆ツール'

# Row 10
$ws.Range("A10").Value = '2025-12-01 01:44:37'
$ws.Range("B10").Value = '【急募】革新的ペット向けECプラットフォーム開発エンジニア募集'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5443928'
$ws.Range("G10").Value = 68
$ws.Range("H10").Value = '◆開発'

# Row 11
$ws.Range("A11").Value = '2025-12-01 01:44:37'
$ws.Range("B11").Value = '【急募】魅力的なWEBサイト制作のフリーランスを探しています!'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5444036'
$ws.Range("G11").Value = 45
$ws.Range("H11").Value = '◇サイト'

# Row 12
$ws.Range("A12").Value = '2025-12-01 01:44:37'
$ws.Range("B12").Value = '【急募】Googledriveのロール管理・共有設定の専門家募集'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5444395'
$ws.Range("G12").Value = 38
$ws.Range("H12").Value = '◇管理'

# Row 13
$ws.Range("A13").Value = '2025-12-01 01:44:37'
$ws.Range("B13").Value = 'wordpressレンダリングを妨げるリソースの除外'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5016989'
$ws.Range("G13").Value = 33
$ws.Range("H13").Value = '○WordPress'

# Row 14
$ws.Range("A14").Value = '2025-12-01 01:44:37'
$ws.Range("B14").Value = 'Salesforce Agentforceの構築・導入支援'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5437485'
$ws.Range("G14").Value = 18
$ws.Range("H14").ClearContents()

# Row 15
$ws.Range("A15").Value = '2025-12-01 01:44:37'
$ws.Range("B15").Value = 'comfyui(paperspace)でエロ動画のループ物を作成したいです。その方法を教えてください'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5444370'
$ws.Range("G15").Value = 10
$ws.Range("H15").ClearContents()

# Row 16
$ws.Range("A16").Value = '2025-12-01 01:44:37'
$ws.Range("B16").Value = '【急募】ミニPCでクラウドストレージ(nextcloud)とOpenWrtルータ化の依頼'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5444262'
$ws.Range("G16").Value = 10
$ws.Range("H16").ClearContents()

# Row 17
$ws.Range("A17").Value = '2025-12-01 01:44:37'
$ws.Range("B17").Value = '空き室情報を拾ってくスクリプト作成'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5444064'
$ws.Range("G17").Value = 10
$ws.Range("H17").ClearContents()

# Re-add hyperlinks F2:F17 in final row order (matches <hyperlinks> block in the diff)
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5443957') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5443998') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5444251') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5444198') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5444378') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5444489') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5251319') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5444141') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5443928') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5444036') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5444395') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5016989') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5437485') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5444370') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5444262') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5444064') | Out-Null

# Hyperlinks.Add() re-applies the Hyperlink cell style via a fresh style index;
# reset explicitly to the workbooks existing "Hyperlink" named style so the F
# column keeps using the original style id (s="1") like the rest of the sheet.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("F12").Style = "Hyperlink"
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("F14").Style = "Hyperlink"
$ws.Range("F15").Style = "Hyperlink"
$ws.Range("F16").Style = "Hyperlink"
$ws.Range("F17").Style = "Hyperlink"